$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE demo.sex IN ['Unknown']
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

$filesQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE demo.sex IN ['Unknown']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis 
'@

$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.sex IN ['Unknown']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Row 2 (CasesTab): B2 -> Cases query (with Cohort line added), C2 -> new stats query
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

# Row 3 (SamplesTab): B3 unchanged, C3 -> new stats query
$ws.Range("C3").Value = $statQuery

# Row 4 (FilesTab): B4 -> Files query (Study Code line removed), C4 -> new stats query
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# Update selection to C2 (matches diff's <selection activeCell="C2" sqref="C2"/>)
$ws.Range("C2").Select() | Out-Null

# Row heights recompute (wrap-text autofit) now that the cell text has changed
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210
